$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ALC_updates = @{
    "H43" = 1663.3334
    "I43" = 1995
    "K43" = 1995
    "M43" = -1926
    "H70" = 1888.7778
    "J70" = 1888.7778
    "L70" = 5666.3334
    "N70" = -6206.3334
    "H73" = 1888.7778
    "J73" = 1888.7778
    "L73" = 5666.3334
    "N73" = -7538.3334
    "H74" = 11979.637
    "I74" = 11577.6
    "K74" = 11577.6
    "M74" = -10641.6
    "H77" = 11979.637
    "I77" = 11577.6
    "K77" = 57888
    "M77" = -53208
    "H98" = 709.4286
    "I98" = 495.6
    "J98" = 1244
    "K98" = 495.6
    "L98" = 1244
    "M98" = 1002.4
    "N98" = -4240
    "H100" = 2772.4666
    "I100" = 2899.0715
    "J100" = 1000
    "K100" = 2899.0715
    "L100" = 1000
    "M100" = -2358.0715
    "N100" = -2082
    "H122" = 709.4286
    "I122" = 495.6
    "J122" = 1244
    "K122" = 1486.8
    "L122" = 3732
    "M122" = 963.1999999999998
    "N122" = -8632
    "H137" = 2652.5715
    "J137" = 3033
    "L137" = 9099
    "N137" = -14199
    "H141" = 0
    "I141" = 0
    "K141" = 0
}
foreach ($cellRef in $ALC_updates.Keys) {
    $ws.Range($cellRef).Value = $ALC_updates[$cellRef]
}
$ws.Range("M141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ARM_updates = @{
    "H23" = 22637.25
    "I23" = 17999.5
    "K23" = 17999.5
    "M23" = -17740.5
    "H31" = 17000
    "I31" = 17000
    "K31" = 17000
    "M31" = -16706
    "H102" = 2757.4443
    "I102" = 2974
    "K102" = 2974
    "M102" = -1352
    "H109" = 47500
    "J109" = 47500
    "L109" = 47500
    "N109" = -50274
}
foreach ($cellRef in $ARM_updates.Keys) {
    $ws.Range($cellRef).Value = $ARM_updates[$cellRef]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$BSM_updates = @{
    "H35" = 60064
    "J35" = 60064
    "L35" = 60064
    "N35" = -60684
    "H80" = 893.8182
    "I80" = 510.4
    "J80" = 1213.3334
    "K80" = 510.4
    "L80" = 1213.3334
    "M80" = 487.6
    "N80" = -3209.3334
    "H83" = 893.8182
    "I83" = 510.4
    "J83" = 1213.3334
    "K83" = 2552
    "L83" = 6066.666999999999
    "M83" = 2440
    "N83" = -16050.667
    "H99" = 26502.25
    "I99" = 51005
    "J99" = 1999.5
    "K99" = 51005
    "L99" = 1999.5
    "M99" = -49507
    "N99" = -4995.5
}
foreach ($cellRef in $BSM_updates.Keys) {
    $ws.Range($cellRef).Value = $BSM_updates[$cellRef]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$CRP_updates = @{
    "H31" = 2445.8333
    "I31" = 1748.3334
    "J31" = 3143.3333
    "K31" = 1748.3334
    "L31" = 3143.3333
    "M31" = -1453.3334
    "N31" = -3733.3333
    "H34" = 2445.8333
    "I34" = 1748.3334
    "J34" = 3143.3333
    "K34" = 1748.3334
    "L34" = 3143.3333
    "M34" = -1546.3334
    "N34" = -3547.3333
    "H50" = 40083.5
    "J50" = 40083.5
    "L50" = 40083.5
    "N50" = -41333.5
    "H60" = 35675.332
    "J60" = 48064
    "L60" = 48064
    "N60" = -49086
    "H62" = 5749.75
    "I62" = 4666.3335
    "J62" = 9000
    "K62" = 4666.3335
    "L62" = 9000
    "M62" = -4042.3335
    "N62" = -10248
    "H65" = 5749.75
    "I65" = 4666.3335
    "J65" = 9000
    "K65" = 23331.6675
    "L65" = 45000
    "M65" = -20211.6675
    "N65" = -51240
    "H105" = 2000
    "I105" = 2000
    "K105" = 2000
    "M105" = -253
    "H107" = 4499.5
    "J107" = 5999
    "L107" = 5999
    "N107" = -9839
    "H134" = 2624.1667
    "J134" = 2749
    "L134" = 8247
    "N134" = -13317
}
foreach ($cellRef in $CRP_updates.Keys) {
    $ws.Range($cellRef).Value = $CRP_updates[$cellRef]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$CUL_updates = @{
    "H94" = 10000
    "J94" = 10000
    "L94" = 30000
    "N94" = -31352
    "H103" = 3695
    "J103" = 5495
    "L103" = 16485
    "N103" = -18243
    "H113" = 1164.6666
    "J113" = 1000
    "L113" = 3000
    "N113" = -7340
}
foreach ($cellRef in $CUL_updates.Keys) {
    $ws.Range($cellRef).Value = $CUL_updates[$cellRef]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$GSM_updates = @{
    "H99" = 5797.6
    "I99" = 5797.6
    "K99" = 5797.6
    "M99" = -3551.6
    "H105" = 9999
    "J105" = 9999
    "L105" = 9999
    "N105" = -16987
}
foreach ($cellRef in $GSM_updates.Keys) {
    $ws.Range($cellRef).Value = $GSM_updates[$cellRef]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$LTW_updates = @{
    "H22" = 6278.143
    "I22" = 3870.5715
    "J22" = 8685.714
    "K22" = 3870.5715
    "L22" = 8685.714
    "M22" = -3575.5715
    "N22" = -9275.714
    "H27" = 6278.143
    "I27" = 3870.5715
    "J27" = 8685.714
    "K27" = 3870.5715
    "L27" = 8685.714
    "M27" = -3763.5715
    "N27" = -8899.714
    "H100" = 1300
    "I100" = 1200
    "K100" = 1200
    "M100" = -659
    "H132" = 3806
    "I132" = 3178.5715
    "K132" = 9535.7145
    "M132" = -7005.7145
}
foreach ($cellRef in $LTW_updates.Keys) {
    $ws.Range($cellRef).Value = $LTW_updates[$cellRef]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$WVR_updates = @{
    "H54" = 34646.668
    "J54" = 34646.668
    "L54" = 34646.668
    "N54" = -35686.668
    "H107" = 743.6667
    "J107" = 874.5
    "L107" = 2623.5
    "N107" = -6463.5
    "H126" = 1509.4814
    "I126" = 1260.5
    "J126" = 1580.619
    "K126" = 3781.5
    "L126" = 4741.857
    "M126" = -1311.5
    "N126" = -9681.857
    "H132" = 9749.5
    "I132" = 4000
    "J132" = 11666
    "K132" = 12000
    "L132" = 34998
    "M132" = -9470
    "N132" = -40058
}
foreach ($cellRef in $WVR_updates.Keys) {
    $ws.Range($cellRef).Value = $WVR_updates[$cellRef]
}
